# Program now automatically tracks all scores from game and writes to
# excel upon exiting.
#
# Each level sheet (s0, s1, s2, s5, s6, s7, s8) logs one row per visit
# ("Visit #" / "Time Spent"), plus a small D column summary showing the
# (repeated) "Mean Time Spent" label next to the mean of the Time Spent
# values recorded so far. This run updates every sheet with the latest
# recorded visits.

$wb = $excel.ActiveWorkbook

# --- s0: single visit recorded (s0_1); second visit row removed ---
$ws = $wb.Worksheets.Item("s0")
$ws.Range("B2").Value = 0.039
$ws.Range("A3:B3").ClearContents()
$ws.Range("D3").Value = 0.039

# --- s1: single visit recorded (s1_1); second visit row removed ---
$ws = $wb.Worksheets.Item("s1")
$ws.Range("B2").Value = 6.767
$ws.Range("A3:B3").ClearContents()
$ws.Range("D3").Value = 6.767

# --- s2: two visits recorded (s2_1, s2_2); third visit row removed ---
$ws = $wb.Worksheets.Item("s2")
$ws.Range("B2").Value = 1.971
$ws.Range("B3").Value = 2.544
$ws.Range("A4:D4").Delete()
$ws.Range("D3").Value = 2.2575

# --- s5: no visits recorded this run; data rows cleared entirely ---
$ws = $wb.Worksheets.Item("s5")
$ws.Range("A2:D3").Delete()

# --- s6: second visit (s6_2) now recorded alongside the first ---
$ws = $wb.Worksheets.Item("s6")
$ws.Range("B2").Value = 6.07
$ws.Range("A3").Value = "s6_2"
$ws.Range("B3").Value = 6.073
$ws.Range("D3").Value = 6.0715

# --- s7: second visit (s7_2) now recorded alongside the first ---
$ws = $wb.Worksheets.Item("s7")
$ws.Range("B2").Value = 61.19
$ws.Range("A3").Value = "s7_2"
$ws.Range("B3").Value = 61.187
$ws.Range("D3").Value = 61.1885

# --- s8: second visit (s8_2) recorded; extra score column G tracked ---
$ws = $wb.Worksheets.Item("s8")
$ws.Range("B2").Value = 3.691
$ws.Range("G2").Value = 2800
$ws.Range("A3").Value = "s8_2"
$ws.Range("B3").Value = 3.655
$ws.Range("D3").Value = 3.673
$ws.Range("G3").Value = 3000
